# Insert a new row at position 52 (shifts existing rows 52-162 down to 53-163)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new weekly record
$ws.Range("A52").Value = 11
$ws.Range("B52").Value = "Vega Monumental Concepción"
$ws.Range("C52").Value = "Bíobío"
$ws.Range("D52").Value = 45044
$ws.Range("E52").Value = 8
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100108
$ws.Range("H52").Value = "Tropicales y subtropicales"
$ws.Range("I52").Value = 100108002
$ws.Range("J52").Value = "Mango"
$ws.Range("K52").Value = "Sin especificar"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 200
$ws.Range("N52").Value = 7500
$ws.Range("O52").Value = 8000
$ws.Range("P52").Value = 7750
$ws.Range("Q52").Value = '$/bandeja 4 kilos'
$ws.Range("R52").Value = "Perú"
$ws.Range("S52").Value = 1938
$ws.Range("T52").Value = 4
